$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "devices"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# New cell E2 = model -> "Galaxy S7"
$ws1.Range("E2").Value = "Galaxy S7"

# Column width tweaks (col A + col E got wider)
$ws1.Columns.Item(1).ColumnWidth = 19
$ws1.Columns.Item(5).ColumnWidth = 26.1875

# Selection / scroll position
$ws1.Activate()
$ws1.Range("E3").Select()

# ---------------------------------------------------------------------------
# Sheet 2: "signIn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Header row renames / additions
$ws2.Range("D1").Value = "Licence_Number"
$ws2.Range("E1").Value = "Licence_StartDate"
$ws2.Range("F1").Value = "Licence_ExpireDate"
$ws2.Range("G1").Value = "class_Type"

# Copy header style (s=1, the yellow-fill header style) onto the new cells
$ws2.Range("C1").Copy() | Out-Null
$ws2.Range("F1:G1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Data row additions
$ws2.Range("D2").Value = "RO4447865"
$ws2.Range("G2").Value = "General"

# Date cells - stored as quoted text with a custom dd-mm-yyyy display format
$ws2.Range("E2").Value = "'09-12-2013"
$ws2.Range("E2").NumberFormat = "dd\-mm\-yyyy"
$ws2.Range("F2").Value = "'08-12-2016"
$ws2.Range("F2").NumberFormat = "dd\-mm\-yyyy"

# Column widths
$ws2.Columns.Item(3).ColumnWidth = 16.45
$ws2.Columns.Item(4).ColumnWidth = 15.17
$ws2.Columns.Item(5).ColumnWidth = 29.02
$ws2.Columns.Item(6).ColumnWidth = 17.74
$ws2.Columns.Item(7).ColumnWidth = 9.59

# Page setup - printed as A4 portrait
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# Selection / scroll position
$ws2.Activate()
$ws2.Range("C1").Select()
